# Update NATMI LR-pair output (Gdf2-Bmpr2) with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value2 = 0.074444
$ws.Range("M2").Value2 = 30.46625333333334
$ws.Range("N2").Value2 = 91.39876000000001
$ws.Range("O2").Value2 = 0.2185380492512374
$ws.Range("P2").Value2 = 0.2331534018544084
$ws.Range("Q2").Value2 = 0.756009921048889
$ws.Range("R2").Value2 = 6.80408928944
$ws.Range("S2").Value2 = 0.2185380492512374
$ws.Range("T2").Value2 = 0.2331534018544084

# Row 3
$ws.Range("H3").Value2 = 0.074444
$ws.Range("O3").Value2 = 0.2491807703757967
$ws.Range("P3").Value2 = 0.2658454419670822
$ws.Range("Q3").Value2 = 0.8620152654613332
$ws.Range("R3").Value2 = 7.758137389151999
$ws.Range("S3").Value2 = 0.2491807703757967
$ws.Range("T3").Value2 = 0.2658454419670822

# Row 4
$ws.Range("H4").Value2 = 0.074444
$ws.Range("M4").Value2 = 23.69037333333334
$ws.Range("N4").Value2 = 71.07112000000001
$ws.Range("O4").Value2 = 0.1699338582153697
$ws.Range("P4").Value2 = 0.181298667526812
$ws.Range("Q4").Value2 = 0.5878687174755556
$ws.Range("R4").Value2 = 5.29081845728
$ws.Range("S4").Value2 = 0.1699338582153697
$ws.Range("T4").Value2 = 0.181298667526812

# Row 5
$ws.Range("H5").Value2 = 0.074444
$ws.Range("M5").Value2 = 26.2168665
$ws.Range("N5").Value2 = 52.433733
$ws.Range("O5").Value2 = 0.18805669340777
$ws.Range("P5").Value2 = 0.1337556791894743
$ws.Range("Q5").Value2 = 0.650562803242
$ws.Range("R5").Value2 = 3.903376819452
$ws.Range("S5").Value2 = 0.18805669340777
$ws.Range("T5").Value2 = 0.1337556791894743

# Row 6
$ws.Range("H6").Value2 = 0.074444
$ws.Range("M6").Value2 = 24.297748
$ws.Range("N6").Value2 = 72.893244
$ws.Range("O6").Value2 = 0.1742906287498262
$ws.Range("P6").Value2 = 0.1859468094622229
$ws.Range("Q6").Value2 = 0.6029405173706666
$ws.Range("R6").Value2 = 5.426464656335999
$ws.Range("S6").Value2 = 0.1742906287498262
$ws.Range("T6").Value2 = 0.1859468094622229
